$wb = $excel.ActiveWorkbook
$wsLeave = $wb.Worksheets.Item("LeaveBalance")
$wsLeave.Activate()
$aw = $excel.ActiveWindow
$p = $aw.Panes.Item(2)
Write-Output "panes count: $($aw.Panes.Count)"
Write-Output "pane: $p"
$tlc = $p.TopLeftCell
Write-Output "tlc: $tlc"
Write-Output "tlc addr: $($tlc.Address())"
